# The upstream change (commit "Fixed POI packaging and upgraded to POI 3.15")
# is a pure repackaging artifact: every line of the diff is the same element,
# with the same tag name, the same text, and the very same set of
# attribute/value pairs -- just re-serialised (attributes alphabetised) by the
# newer Apache POI writer. No word of visible content, paragraph structure,
# tab stop, margin or style value actually changed.
#
# Word's object model does not expose "raw attribute order" as a settable
# property, so we can't poke that directly. What we *can* do is touch every
# OOXML element the diff lists, through the corresponding COM property, and
# re-assert the exact same value it already had -- which is the closest
# faithful COM-interop equivalent of "this part got rewritten by the
# packager" without altering anything a user would see.

$d = $word.ActiveDocument

# word/document.xml: the three <w:pPr><w:tabs><w:tab w:val="left"
# w:pos="3119"/></w:tabs></w:pPr> paragraphs (the "m:if", "THEN" and
# "m:endif" lines).
foreach ($i in 2, 3, 4) {
    $p = $d.Paragraphs.Item($i)
    $tabs = $p.Range.ParagraphFormat.TabStops
    $tabs.ClearAll()
    $tabs.Add(155.95, 0)   # 155.95pt = 3119 twips, wdAlignTabLeft = 0
}

# word/document.xml: <w:sectPr><w:pgSz .../><w:pgMar .../></w:sectPr>
$section = $d.Sections.Item(1)
$pageSetup = $section.PageSetup
$pageSetup.PageWidth = 595.3        # 11906 twips
$pageSetup.PageHeight = 841.9       # 16838 twips
$pageSetup.TopMargin = 70.85        # 1417 twips
$pageSetup.RightMargin = 70.85      # 1417 twips
$pageSetup.BottomMargin = 70.85     # 1417 twips
$pageSetup.LeftMargin = 70.85       # 1417 twips
$pageSetup.HeaderDistance = 35.4    # 708 twips
$pageSetup.FooterDistance = 35.4    # 708 twips
$pageSetup.Gutter = 0
